$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H matches the width of the other data columns (C:G), ~12.71
$ws.Range("H1").ColumnWidth = 11.8

# Header row
$ws.Range("H1").Value = "NRMSE"

# Updated data values (rows 2-6, columns A-G) plus new column H
$data = @(
    @(3.2274263923235722, 2.4330235866390177, 0.43516490023005738, 0.65967029661040322, 0.67923606550503024, 0.56483509976994262, 0.77120646395525005, 0.23489275053301109),
    @(4.1396993045858537, 3.2055406142659377, 0.55488456300778066, 0.74490574102216489, 0.75799021382500298, 0.44511543699221934, 0.69829100786795406, 0.28888341274151103),
    @(3.1179784642106232, 2.3569886736135399, 0.31426912593927947, 0.56059711552886127, 0.5291847044484822,  0.68573087406072053, 0.83136303180551141, 0.21239635314786259),
    @(3.470872981519741,  2.571916864747716,  0.36996448377599483, 0.60824705817290625, 0.59721743056953813, 0.63003551622400522, 0.7954007331535532,  0.24783098761297689),
    @(3.7358730822107589, 2.8809497485340585, 0.52984077104432237, 0.72790162181734586, 0.76866321999307863, 0.47015922895567763, 0.72216651386117059, 0.27755372081803564)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
